$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data (2021-12-22) is published, so it gets inserted right
# after the existing header + first block of rows, pushing every later row
# down by 4 (one new row per "Calidad": Especial, Primera, Segunda, Tercera).
$ws.Range("A6:A9").EntireRow.Insert()

# Constant columns (A,B,C,E,F,G,H,I,J,K,R) are identical for every record in
# this sheet - reuse the same values used throughout the table.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$prodId    = 100108
$producto  = "Tropicales y subtropicales"
$catId     = 100108005
$categoria = "Piña"
$variedad  = "Caramelo"
$origen    = "Ecuador"
$fecha     = "2021-12-22"

function Set-Row($r, $calidad, $volumen, $pmin, $pmax, $pprom, $unidad, $precioKg, $kgUnidad) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $prodId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $catId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $pmin
    $ws.Cells.Item($r, 15).Value = $pmax
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 6 "Especial" 200 16000 17000 16500 "$/caja 10 unidades" 1650 10
Set-Row 7 "Primera"  250 16000 17000 16500 "$/caja 12 unidades" 1375 12
Set-Row 8 "Segunda"  260 16000 17000 16500 "$/caja 14 unidades" 1179 14
Set-Row 9 "Tercera"  270 16000 17000 16500 "$/caja 16 unidades" 1031 16
